# Update row power flow results (pl_mw) for the 380 kV case.
# Columns B,C and E..I on rows 2..25 receive new computed values;
# columns A, D and J..O are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bcData = @(
        @(1.196066065259174, 0.3842460153423417),
        @(1.045931578302486, 0.3359693454906676),
        @(0.9535825454929636, 0.3062284184046575),
        @(0.9159082741816178, 0.2940835692519101),
        @(0.9096500097268745, 0.2920653965045688),
        @(0.9530746236832215, 0.3060647310758213),
        @(1.144334302877667, 0.367620524531219),
        @(1.518081104647422, 0.487567922026642),
        @(1.791908278422341, 0.5752665516382081),
        @(1.916323874004888, 0.6150794675478437),
        @(1.963415610197046, 0.6301444320602059),
        @(1.953274532292653, 0.6269004191308909),
        @(1.920198576381438, 0.6163190957065012),
        @(1.899935760603626, 0.6098362677519731),
        @(1.783774393338717, 0.5726630763650178),
        @(1.712474627488916, 0.5498379087127887),
        @(1.671450800773243, 0.5367017662848639),
        @(1.657558426021581, 0.5322527577242226),
        @(1.720066068981964, 0.552268481733563),
        @(1.929914378081833, 0.6194273907634624),
        @(2.066935550713822, 0.6632538316804357),
        @(1.993816431376729, 0.639868711851534),
        @(1.716634079949586, 0.551169662230393),
        @(1.417109775708695, 0.4551967779846677)
    )

$eiData = @(
        @(0.4267877700075502, 0.4443680307746121, 0.2988891640978295, 0.4179759268276513, 0.2534748189175851),
        @(0.372013532808964, 0.3878228170618172, 0.2908448685241751, 0.4214889163782658, 0.2609967475567174),
        @(0.3385037565169284, 0.3531389305169483, 0.2865333117792801, 0.4241743304012005, 0.2660843967949091),
        @(0.3248762196165416, 0.3390132514313251, 0.2849316345146065, 0.4254004957568753, 0.2682745533835309),
        @(0.322614989461286, 0.336668177824194, 0.2846749845204926, 0.4256120319040235, 0.2686452561935369),
        @(0.3383198610164726, 0.3529483938344953, 0.2865110854549187, 0.4241903344754832, 0.2661134621239576),
        @(0.4078746997522842, 0.4248636149813478, 0.2959836977471753, 0.4190769038325044, 0.2559703601887477),
        @(0.5453655366994923, 0.5661985755041457, 0.3196545790892173, 0.4132895915187191, 0.2398504570444047),
        @(0.6472337432702773, 0.6702781546542269, 0.3403125575931369, 0.4116886565497708, 0.2303721130619678),
        @(0.6938031159386213, 0.7176906081379002, 0.3504547812095211, 0.4115503974662289, 0.226588039633004),
        @(0.7114739079612349, 0.7356546913071611, 0.3544053144905632, 0.4115840669493309, 0.2252321834302435),
        @(0.7076665373244424, 0.7317853510981394, 0.3535495642509829, 0.4115729706922764, 0.2255207433185191),
        @(0.6952561638514396, 0.7191683204515869, 0.3507775754691806, 0.4115514365946211, 0.2264749405507303),
        @(0.6876592290190757, 0.7114413442032514, 0.349094046863982, 0.4115494846165006, 0.2270694894522585),
        @(0.6441951930320613, 0.6671810134426437, 0.3396649649268824, 0.4117096654860291, 0.2306301326450537),
        @(0.6175922234768052, 0.6400460337215605, 0.3340732703470053, 0.4119598856088231, 0.2329504807406195),
        @(0.6023122693328418, 0.6244449056556647, 0.3309270120017658, 0.4121592692553691, 0.2343346776852009),
        @(0.5971423177658011, 0.6191636801734006, 0.3298736699415059, 0.412236267221715, 0.2348118240449111),
        @(0.6204219212730493, 0.642933953830422, 0.3346612574804624, 0.4119275004682805, 0.2326983341953834),
        @(0.6989003902806843, 0.7228739723492197, 0.3515887709426551, 0.4115554175086231, 0.2261925677832366),
        @(0.7504017872733471, 0.7751780083420101, 0.3632937329212353, 0.4118142167468619, 0.2223905834395659),
        @(0.722894177605113, 0.7472568307915566, 0.3569869168392046, 0.4116297664509148, 0.2243782044791693),
        @(0.6191425701634472, 0.6416283278902313, 0.3343952153997662, 0.4119419689084083, 0.2328121733759332),
        @(0.5080343731648469, 0.5279251897347308, 0.31268895084132, 0.4143944794371635, 0.2438008223960892)
    )

$bcArray = New-Object 'object[,]' $bcData.Count, 2
for ($i = 0; $i -lt $bcData.Count; $i++) {
    for ($j = 0; $j -lt 2; $j++) {
        $bcArray[$i, $j] = $bcData[$i][$j]
    }
}

$eiArray = New-Object 'object[,]' $eiData.Count, 5
for ($i = 0; $i -lt $eiData.Count; $i++) {
    for ($j = 0; $j -lt 5; $j++) {
        $eiArray[$i, $j] = $eiData[$i][$j]
    }
}

$ws.Range("B2:C25").Value2 = $bcArray
$ws.Range("E2:I25").Value2 = $eiArray
